# Add a "Team_Division" column (G) to Sheet1, filling each team's
# NFL division, and update the sheet view's active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map each data row (2-33) to its division name.
$divisions = @{
    2  = "NFC West"
    3  = "NFC South"
    4  = "AFC North"
    5  = "AFC East"
    6  = "NFC South"
    7  = "NFC North"
    8  = "AFC North"
    9  = "AFC North"
    10 = "NFC East"
    11 = "AFC West"
    12 = "NFC North"
    13 = "NFC North"
    14 = "AFC South"
    15 = "AFC South"
    16 = "AFC South"
    17 = "AFC West"
    18 = "AFC West"
    19 = "NFC West"
    20 = "AFC East"
    21 = "NFC North"
    22 = "AFC East"
    23 = "NFC South"
    24 = "NFC East"
    25 = "AFC East"
    26 = "AFC West"
    27 = "NFC East"
    28 = "AFC North"
    29 = "NFC West"
    30 = "NFC West"
    31 = "NFC South"
    32 = "AFC South"
    33 = "NFC East"
}

# Seed the shared-string table so brand-new unique values are created in
# the same order the original workbook used (NFC West, AFC West, NFC
# South, AFC North, AFC East, NFC North, NFC East, AFC South), by writing
# to the first row that introduces each distinct division first.
$seedRows = @(2, 11, 3, 4, 5, 7, 10, 14)
foreach ($r in $seedRows) {
    $ws.Cells.Item($r, 7).Value = $divisions[$r]
}

# Now fill in the Team_Division column for every remaining row.
foreach ($r in 2..33) {
    $ws.Cells.Item($r, 7).Value = $divisions[$r]
}

# Update the header cell for column G (already set by the author, but make
# sure it is present/correct).
$ws.Cells.Item(1, 7).Value = "Team_Division"

# Update the view: active cell E24 with the window scrolled so row 14 /
# column C is toward the top-left.
$ws.Range("E24").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$win.ScrollColumn = 3
